# Adds 6 new collected-data rows (rows 40-45) to the "coletaDeDados" sheet,
# matching the rows already logged by the translation widget.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: spoken language, translation language, sent text, translated text,
# browser, mobile device flag, OS, computer flag, tablet flag, accessed time.
$rows = @(
    @("Português", "Inglês",    "Bom dia, gostaria de falar com o senhor pedro", "Good morning, i would like to talk to mr. pedro", "Chrome", $false, "Windows", $true, $false, "23-06-2024 11:25:11"),
    @("Português", "Inglês",    "Seria possível?",                                "It would be possible?",                           "Chrome", $false, "Windows", $true, $false, "23-06-2024 11:25:32"),
    @("Inglês",    "Português", "Yes, one momento please!",                       "Sim, um momento, por favor!",                     "Chrome", $false, "Windows", $true, $false, "23-06-2024 11:25:47"),
    @("Português", "Inglês",    "Posso me sentar aqui?",                          "May i sit here?",                                 "Chrome", $false, "Windows", $true, $false, "23-06-2024 11:26:01"),
    @("Inglês",    "Português", "Yes, feel free",                                 "Sim, sinta -se livre",                            "Chrome", $false, "Windows", $true, $false, "23-06-2024 11:26:26"),
    @("Inglês",    "Português", "Mr. pedro will talk to you soon, okay?",         "O sr. pedro falará com você em breve, ok?",       "Chrome", $false, "Windows", $true, $false, "23-06-2024 11:27:04")
)

$startRow = 40
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value  = $data[0]   # A - Idioma de fala
    $ws.Cells.Item($r, 2).Value  = $data[1]   # B - Idioma de Tradução
    $ws.Cells.Item($r, 3).Value  = $data[2]   # C - Texto enviado
    $ws.Cells.Item($r, 4).Value  = $data[3]   # D - Texto traduzido
    $ws.Cells.Item($r, 5).Value  = $data[4]   # E - Navegador
    $ws.Cells.Item($r, 6).Value  = $data[5]   # F - Dispositivo móvel
    $ws.Cells.Item($r, 7).Value  = $data[6]   # G - Sistema Operacional
    $ws.Cells.Item($r, 8).Value  = $data[7]   # H - Computador
    $ws.Cells.Item($r, 9).Value  = $data[8]   # I - Tablet
    $ws.Cells.Item($r, 10).Value = $data[9]   # J - Horário acessado
}
